$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue date / volume number) ---
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Data table updates (rows 15-29) ---
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -66.666666666666
$ws.Range("N15").Value = -34.883720930232
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -13.043478260869
$ws.Range("I16").Value = 232
$ws.Range("J16").Value = 215
$ws.Range("K16").Value = 7.906976744186
$ws.Range("L16").Value = 34.883720930232
$ws.Range("M16").Value = -15.636363636363
$ws.Range("N16").Value = -84.206943498978
# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -64.285714285714
$ws.Range("F17").Value = 42
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = 2.439024390243
$ws.Range("I17").Value = 432
$ws.Range("J17").Value = 425
$ws.Range("K17").Value = 1.647058823529
$ws.Range("L17").Value = 24.495677233429
$ws.Range("M17").Value = 82.278481012658
$ws.Range("N17").Value = -34.045801526717
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 185
$ws.Range("J18").Value = 232
$ws.Range("K18").Value = -20.258620689655
$ws.Range("L18").Value = -8.866995073891
$ws.Range("M18").Value = -55.314009661835
$ws.Range("N18").Value = -90.439276485788
# Row 19
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = 4.225352112676
$ws.Range("I19").Value = 773
$ws.Range("J19").Value = 745
$ws.Range("K19").Value = 3.758389261744
$ws.Range("L19").Value = 33.275862068965
$ws.Range("M19").Value = 39.028776978417
$ws.Range("N19").Value = -18.545837723919
# Row 20
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = 25.806451612903
$ws.Range("I20").Value = 378
$ws.Range("J20").Value = 304
$ws.Range("K20").Value = 24.342105263157
$ws.Range("L20").Value = 79.146919431279
$ws.Range("M20").Value = 61.538461538461
$ws.Range("N20").Value = -83.636363636363
# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -28.888888888888
$ws.Range("F21").Value = 189
$ws.Range("G21").Value = 184
$ws.Range("H21").Value = 2.717391304347
$ws.Range("I21").Value = 2033
$ws.Range("J21").Value = 1949
$ws.Range("K21").Value = 4.309902514109
$ws.Range("L21").Value = 31.245965138799
$ws.Range("M21").Value = 17.175792507204
$ws.Range("N21").Value = -72.471225457007
# Row 22
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 66.666666666666
$ws.Range("M22").Value = 29.629629629629
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -21.052631578947
$ws.Range("I23").Value = 210
$ws.Range("J23").Value = 184
$ws.Range("K23").Value = 14.130434782608
$ws.Range("L23").Value = 18.644067796610
$ws.Range("M23").Value = 59.090909090909
# Row 24
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 34.375
$ws.Range("F24").Value = 161
$ws.Range("G24").Value = 180
$ws.Range("H24").Value = -10.555555555555
$ws.Range("I24").Value = 1993
$ws.Range("J24").Value = 2284
$ws.Range("K24").Value = -12.740805604203
$ws.Range("L24").Value = 1.787538304392
$ws.Range("M24").Value = 54.856254856254
# Row 25
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 82
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = 10.810810810810
$ws.Range("I25").Value = 740
$ws.Range("J25").Value = 793
$ws.Range("K25").Value = -6.683480453972
$ws.Range("L25").Value = -4.392764857881
$ws.Range("M25").Value = -1.726427622841
# Row 26
$ws.Range("D26").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 36
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = 2.272727272727
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -54.545454545454
$ws.Range("I27").Value = 95
$ws.Range("J27").Value = 81
$ws.Range("K27").Value = 17.283950617283
$ws.Range("L27").Value = -1.041666666666
# Row 28
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = 26.315789473684
# Row 29
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 17
$ws.Range("K29").Value = 29.411764705882

$excel.CutCopyMode = $false
Write-Host "Edits applied"
